$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell. Values that look like plain
# numbers (e.g. "01", "0230", "120") must be forced to Text so Excel
# doesn't silently coerce them into numbers and eat leading zeros - the
# source workbook stores all of these as shared strings. We do this with
# the classic leading-apostrophe ("quote prefix") trick and then restore
# the "Normal" cell style so the cell ends up with the same (default)
# formatting it would have had if the text had never looked numeric.
function Set-TextCell {
    param($range, [string]$text)
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

# ---- Row 1 (helper/index row, B1:M1 = 0..11), bold/boxed style ----
# B1:K1 already correct (0..9); just extend with the two new columns.
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:M1").PasteSpecial(-4122) | Out-Null

# ---- Column A (row index 0..5), bold/boxed style ----
$ws.Range("A7").Value = 5
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

# ---- Clear the two cells that disappear from row 6 (day/start-time
#      columns no longer apply to the "Private Reading" rows) ----
$ws.Range("G6").ClearContents() | Out-Null
$ws.Range("I6").ClearContents() | Out-Null

# ---- Row 2: ESOL Level 2 MWF, 1000-1050am ----
Set-TextCell $ws.Range("B2") "120"
Set-TextCell $ws.Range("C2") "01"
Set-TextCell $ws.Range("D2") "17630"
Set-TextCell $ws.Range("E2") "4"
Set-TextCell $ws.Range("F2") "F"
Set-TextCell $ws.Range("G2") "MWF"
Set-TextCell $ws.Range("H2") "ESOL Level 2"
Set-TextCell $ws.Range("I2") "1000"
Set-TextCell $ws.Range("J2") "1050am"
Set-TextCell $ws.Range("K2") "SEVE"
Set-TextCell $ws.Range("L2") "204"
Set-TextCell $ws.Range("M2") "Fekete Larissa"

# ---- Row 3: ESOL Level 3 MWF, 0230-0320pm ----
Set-TextCell $ws.Range("B3") "130"
Set-TextCell $ws.Range("C3") "01"
Set-TextCell $ws.Range("D3") "17631"
Set-TextCell $ws.Range("E3") "4"
Set-TextCell $ws.Range("F3") "F"
Set-TextCell $ws.Range("G3") "MWF"
Set-TextCell $ws.Range("H3") "ESOL Level 3"
Set-TextCell $ws.Range("I3") "0230"
Set-TextCell $ws.Range("J3") "0320pm"
Set-TextCell $ws.Range("K3") "CBIB"
Set-TextCell $ws.Range("L3") "216"
Set-TextCell $ws.Range("M3") "Fekete Larissa"

# ---- Row 4: Expository Writing for ESOL MWF, section 01, 1100-1150am ----
Set-TextCell $ws.Range("B4") "140"
Set-TextCell $ws.Range("C4") "01"
Set-TextCell $ws.Range("D4") "19338"
Set-TextCell $ws.Range("E4") "4"
Set-TextCell $ws.Range("F4") "F"
Set-TextCell $ws.Range("G4") "MWF"
Set-TextCell $ws.Range("H4") "Expository Writing for ESOL"
Set-TextCell $ws.Range("I4") "1100"
Set-TextCell $ws.Range("J4") "1150am"
Set-TextCell $ws.Range("K4") "SEVE"
Set-TextCell $ws.Range("L4") "204"
Set-TextCell $ws.Range("M4") "Fekete Larissa"

# ---- Row 5: NEW section - Expository Writing for ESOL, section 02, 0330-0420pm ----
Set-TextCell $ws.Range("B5") "140"
Set-TextCell $ws.Range("C5") "02"
Set-TextCell $ws.Range("D5") "19339"
Set-TextCell $ws.Range("E5") "4"
Set-TextCell $ws.Range("F5") "F"
Set-TextCell $ws.Range("G5") "MWF"
Set-TextCell $ws.Range("H5") "Expository Writing for ESOL"
Set-TextCell $ws.Range("I5") "0330"
Set-TextCell $ws.Range("J5") "0420pm"
Set-TextCell $ws.Range("K5") "CBIB"
Set-TextCell $ws.Range("L5") "216"
Set-TextCell $ws.Range("M5") "Fekete Larissa"

# ---- Row 6: Private Reading - Full (995F) ----
Set-TextCell $ws.Range("B6") "995F"
Set-TextCell $ws.Range("C6") "01"
Set-TextCell $ws.Range("D6") "19289"
Set-TextCell $ws.Range("E6") "4"
Set-TextCell $ws.Range("F6") "F"
Set-TextCell $ws.Range("H6") "Private Reading"
Set-TextCell $ws.Range("J6") "Full"
Set-TextCell $ws.Range("K6") "TBA"
Set-TextCell $ws.Range("L6") "TBA"
Set-TextCell $ws.Range("M6") "Fekete Larissa"

# ---- Row 7 (new row): Private Reading - Half (995H) ----
Set-TextCell $ws.Range("B7") "995H"
Set-TextCell $ws.Range("C7") "01"
Set-TextCell $ws.Range("D7") "19290"
Set-TextCell $ws.Range("E7") "2"
Set-TextCell $ws.Range("F7") "F"
Set-TextCell $ws.Range("H7") "Private Reading"
Set-TextCell $ws.Range("J7") "Half"
Set-TextCell $ws.Range("K7") "TBA"
Set-TextCell $ws.Range("L7") "TBA"
Set-TextCell $ws.Range("M7") "Fekete Larissa"

Write-Output "Edit complete"
